# Task: Completed daily operations, 8 hours, 11/09
# Append a new time-log entry (row 20) to Sheet1, continuing the existing
# "Internship" / daily-operations-description pattern used by the prior rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row: Date | Name of Task | Description
$ws.Range("A20").Value = 45239
$ws.Range("A20").NumberFormat = $ws.Range("A19").NumberFormat()
$ws.Range("B20").Value = $ws.Range("B19").Value()
$ws.Range("C20").Value = $ws.Range("C19").Value()

# Match the saved selection state after data entry (cursor moved to next row)
$ws.Range("C21").Select()
